$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("T1").Value = "city/state"
$ws.Range("U1").Value = "exp_per_student"

# Copy the header style (bold, centered, bordered) from the existing header row
# onto the two new header cells so they match the rest of row 1.
$ws.Range("S1").Copy()
$ws.Range("T1:U1").PasteSpecial(-4122)

# city/state = city_location & ", " & state_location for every data row
$cityState = @(
    'Nashville, TN',
    'Alamo, TN',
    'Alcoa, TN',
    'Jamestown, TN',
    'Clinton, TN',
    'Arlington, TN',
    'Athens, TN',
    'Bartlett, TN',
    'Shelbyville, TN',
    'Bells, TN',
    'Camden, TN',
    'Pikeville, TN',
    'Maryville, TN',
    'Bradford, TN',
    'Cleveland, TN',
    'Bristol, TN',
    'Jacksboro, TN',
    'Woodbury, TN',
    'Huntingdon, TN',
    'Elizabethton, TN',
    'Ashland City, TN',
    'Henderson, TN',
    'Tazewell, TN',
    'Celina, TN',
    'Cleveland, TN',
    'Clinton, TN',
    'Newport, TN',
    'Manchester, TN',
    'Collierville, TN',
    'Alamo, TN',
    'Crossville, TN',
    'Nashville, TN',
    'Dayton, TN',
    'Decaturville, TN',
    'Smithville, TN',
    'Nashville, TN',
    'Dickson, TN',
    'Dyersburg, TN',
    'Dyersburg, TN',
    'Elizabethton, TN',
    'Etowah, TN',
    'Somerville, TN',
    'Fayetteville, TN',
    'Jamestown, TN',
    'Winchester, TN',
    'Franklin, TN',
    'Germantown, TN',
    'Dyer, TN',
    'Pulaski, TN',
    'Rutledge, TN',
    'Greeneville, TN',
    'Greeneville, TN',
    'Altamont, TN',
    'Morristown, TN',
    'Chattanooga, TN',
    'Sneedville, TN',
    'Bolivar, TN',
    'Savannah, TN',
    'Rogersville, TN',
    'Brownsville, TN',
    'Lexington, TN',
    'Paris, TN',
    'Centerville, TN',
    'Bruceton, TN',
    'Erin, TN',
    'Humboldt, TN',
    'Waverly, TN',
    'Huntingdon, TN',
    'Gainesboro, TN',
    'Dandridge, TN',
    'Johnson City, TN',
    'Mountain City, TN',
    'Kingsport, TN',
    'Knoxville, TN',
    'Tiptonville, TN',
    'Lakeland, TN',
    'Ripley, TN',
    'Lawrenceburg, TN',
    'Lebanon, TN',
    'Lenoir City, TN',
    'Hohenwald, TN',
    'Lexington, TN',
    'Fayetteville, TN',
    'Loudon, TN',
    'Lafayette, TN',
    'Jackson, TN',
    'Manchester, TN',
    'Jasper, TN',
    'Lewisburg, TN',
    'Maryville, TN',
    'Columbia, TN',
    'McKenzie, TN',
    'Athens, TN',
    'Selmer, TN',
    'Decatur, TN',
    'Memphis, TN',
    'Milan, TN',
    'Millington, TN',
    'Madisonville, TN',
    'Clarksville, TN',
    'Lynchburg, TN',
    'Wartburg, TN',
    'Murfreesboro, TN',
    'Newport, TN',
    'Oak Ridge, TN',
    'Union City, TN',
    'Oneida, TN',
    'Livingston, TN',
    'Paris, TN',
    'Linden, TN',
    'Byrdstown, TN',
    'Benton, TN',
    'Cookeville, TN',
    'Dayton, TN',
    'South Pittsburg, TN',
    'Kingston, TN',
    'Springfield, TN',
    'Rogersville, TN',
    'Murfreesboro, TN',
    'Huntsville, TN',
    'Dunlap, TN',
    'Sevierville, TN',
    'Carthage, TN',
    'Clarksburg, TN',
    'Dover, TN',
    'Blountville, TN',
    'Gallatin, TN',
    'Sweetwater, TN',
    'Nashville, TN',
    'Nashville, TN',
    'Knoxville, TN',
    'Covington, TN',
    'Trenton, TN',
    'Hartsville, TN',
    'Tullahoma, TN',
    'Erwin, TN',
    'Union City, TN',
    'Maynardville, TN',
    'Spencer, TN',
    'McMinnville, TN',
    'Jonesborough, TN',
    'Waynesboro, TN',
    'Dresden, TN',
    'Atwood, TN',
    'Jackson, TN',
    'Sparta, TN',
    'Franklin, TN',
    'Lebanon, TN'
)

# exp_per_student = exp_total / enrollment for every data row
# ($null where exp_total is blank, same rows that are blank in H:S)
$expPerStudent = @(
    19316.2923231646,
    9494.809688581316,
    11587.516960651288,
    $null,
    11690.675953775526,
    10295.232146575057,
    10234.739178690344,
    13881.287726358149,
    10234.884492096828,
    9619.047619047618,
    10666.666666666666,
    11609.431680773881,
    10765.389876880985,
    10128.455284552845,
    9240.919271594119,
    11053.363117398858,
    10171.326029467322,
    10413.941018766756,
    1700000.0,
    11268.052057094877,
    10216.01234356249,
    9787.620064034152,
    10324.078624078624,
    10570.647219690063,
    11259.188626907073,
    9934.917355371901,
    10513.482369209496,
    10799.317406143346,
    11844.799654390323,
    9507.48502994012,
    9306.541019955654,
    14857.125144139563,
    9713.443396226416,
    10301.181102362205,
    8852.280462899931,
    $null,
    9875.753475212203,
    10886.145404663923,
    11401.45426712591,
    10092.768444119796,
    10274.05247813411,
    10643.298332303892,
    11798.31223628692,
    9446.69603524229,
    10948.641087130296,
    19254.545454545456,
    14181.048788582808,
    9109.793033821303,
    10205.17711171662,
    11323.67149758454,
    10380.526735833999,
    11125.792459125792,
    11361.654555617664,
    10058.120322674702,
    10453.422508823145,
    11327.044025157233,
    12308.223289315727,
    10778.041543026706,
    11189.403553299493,
    11068.724733553841,
    9444.745929945733,
    10746.84804246848,
    10755.590551181102,
    9325.227963525836,
    10471.506635441061,
    12778.378378378378,
    10383.587077608616,
    9847.940074906366,
    10634.996582365004,
    10087.062187276626,
    10494.313210848644,
    4570.911285455642,
    10683.967112024666,
    9621.524856044482,
    14710.280373831776,
    11061.444652908067,
    12655.142103721066,
    9458.902161547212,
    9670.69414830609,
    10211.059190031152,
    9506.478209658422,
    10395.927601809955,
    9736.631684157921,
    9917.45330387628,
    9007.090522335146,
    11348.631950573697,
    11085.91282375237,
    9828.787878787878,
    9836.662967073622,
    11460.412625639217,
    9406.25,
    9072.2049689441,
    9889.56937799043,
    10505.873140172278,
    9952.598515134208,
    12974.871808883667,
    10189.879759519037,
    12653.727901614144,
    10627.513438184353,
    9484.280298145617,
    17865.704772475026,
    10379.56204379562,
    9981.079931972788,
    10499.28673323823,
    13101.39084861923,
    10655.637648570511,
    8509.227614490772,
    9948.09800128949,
    11311.67192429022,
    11379.79420018709,
    11068.32298136646,
    10704.37616387337,
    9605.884308510638,
    10488.633585920314,
    14208.791208791208,
    10490.384615384615,
    10001.477746870654,
    11073.71794871795,
    10966.257124262922,
    11211.965134706814,
    10265.878877400295,
    11598.616395710827,
    10213.982358706306,
    10436.923076923076,
    10020.479520479521,
    11828.997613365154,
    9371.534556813745,
    8732.166890982504,
    $null,
    $null,
    $null,
    9441.896770416626,
    12045.790251107828,
    9604.122245913291,
    10798.132381213953,
    10745.045045045044,
    11048.051948051949,
    6104.394679238929,
    12086.901763224181,
    10280.222325150533,
    9638.32870426277,
    10841.8131359852,
    9714.32154502132,
    10745.519713261649,
    $null,
    13715.392561983472,
    11580.256574423182,
    10071.154224935986
)

for ($i = 0; $i -lt $cityState.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 20).Value = $cityState[$i]
    if ($expPerStudent[$i] -ne $null) {
        $ws.Cells.Item($row, 21).Value = $expPerStudent[$i]
    }
}
